$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the "Förändrad" (Changed) date column (C) for rows 2-11
# from 45174 (2023-09-05) to 45175 (2023-09-06)
for ($row = 2; $row -le 11; $row++) {
    $cell = $ws.Range("C$row")
    if ($cell.Value2 -eq 45174) {
        $cell.Value2 = 45175
    }
}
